$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Fechas de la campaña para constelación de Géminis 2022: 14-23 de febrero, 14-24 de marzo",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "2022 Fechas de la campaña para constelación de Géminis: 14-23 de febrero, 14-24 de marzo",
    2
)
